# Small RGB() helper -- PowerShell has no built-in RGB(), and Excel's COM
# Font/Interior .Color property wants a single BGR-packed long.
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Abass Yekeen's "P" on C2: it was left in the plain/unstyled look while
# every other "Present" cell in the row uses the green font style used by B2.
$ws.Range("C2").Font.Color = (RGB 0 176 80)

# Forhad Hussain (row 5) was marked Present for the last three weeks
# (J5:L5); update the attendance to Absent (red) for those dates.
$ws.Range("J5:L5").Font.Color = (RGB 255 0 0)
$ws.Range("J5:L5").Value = "A"

# Leave the cursor/selection where the user last left it before saving.
$ws.Range("B2:C2").Select()
